$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (values + % changes), preserving text type
# and default (unstyled) cell formatting for every touched cell.

$ws.Range('D2').Value = '''26.625.55'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  +0.00%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.597.03'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +0.34%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.10%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''211.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -0.23%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = '''  +0.54%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.06%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E9').Value = '''  -0.02%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D11').Value = '''0.0836'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +0.22%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''1.821.66'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  +0.37%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''1.618.81'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +1.62%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = '''  -0.23%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = '''  -0.42%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''64.81'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -0.23%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''26.627.77'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E18').Value = '''  +0.80%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = '''  +0.10%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''208.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  +0.14%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''7.05'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  +5.40%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''4.27'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  +0.29%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''2.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  -0.17%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = '''  -0.03%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''145.44'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -0.41%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = '''  +0.08%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = '''  +0.23%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = '''  +0.56%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''15.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -0.46%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''0.0507'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +0.30%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = '''  +0.20%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = '''  -0.21%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = '''  +0.77%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''0.624'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -5.51%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''1.273.25'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -1.68%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = '''  +0.32%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = '''  +0.01%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = '''  -0.69%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = '''  +1.15%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''5.48'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  +2.60%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = '''  +1.25%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = '''WEMIXToken'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = '''https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = '''0.957'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +18.39%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = '''TrustWalletToken'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = '''0.785'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -0.69%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = '''Aave'
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = '''64.04'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  +1.46%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''1.733.73'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +0.37%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''89.94'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  +1.04%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = '''  +0.19%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = '''  +4.28%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''0.0509'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  +1.12%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = '''EnergySwap'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = '''7.50'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  +0.13%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = '''USDD'
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = '''https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = '''1.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +0.08%  '
$ws.Range('E51').Style = 'Normal'
